$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first 9 digits of the CPF (row 2, columns A-I)
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 1

# Update the active selection as recorded in the saved view
$ws.Range("L8").Select()
